$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.836.27"
$ws.Range("D3").Value = "2.896.62"
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.09"
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.32"
$ws.Range("D6").Style = $origStyle
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("D9").Style = $origStyle
$ws.Range("D10").Value = "2.894.00"
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.418"
$ws.Range("D11").Style = $origStyle
$ws.Range("D14").Value = "3.422.54"
$ws.Range("D15").Value = "75.677.85"
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("D16").Style = $origStyle
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.36"
$ws.Range("D17").Style = $origStyle
$ws.Range("D18").Value = "2.892.73"
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.88"
$ws.Range("D19").Style = $origStyle
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.59"
$ws.Range("D20").Style = $origStyle
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.27"
$ws.Range("D21").Style = $origStyle
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("D22").Style = $origStyle
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.15"
$ws.Range("D23").Style = $origStyle
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = $origStyle
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.98"
$ws.Range("D25").Style = $origStyle
$ws.Range("D26").Value = "3.045.14"
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.21"
$ws.Range("D27").Style = $origStyle
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("D28").Style = $origStyle
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.987"
$ws.Range("D30").Style = $origStyle
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "501.98"
$ws.Range("D32").Style = $origStyle
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.70"
$ws.Range("D33").Style = $origStyle
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = $origStyle
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.41"
$ws.Range("D36").Style = $origStyle
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.98"
$ws.Range("D37").Style = $origStyle
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.40"
$ws.Range("D41").Style = $origStyle
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.99"
$ws.Range("D43").Style = $origStyle
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.99"
$ws.Range("D47").Style = $origStyle
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.575"
$ws.Range("D49").Style = $origStyle
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.650"
$ws.Range("D51").Style = $origStyle

$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  +11.36%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("E32").Value = "  -6.36%  "
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  -7.61%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  -5.97%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("E51").Value = "  +5.05%  "
